# TornadoEngine plan sheet: collapse the empty spacer row (old r=2), which
# pulls the two task rows up (old r=3,4 -> r=2,3) and the trailing blank
# spacer row up too (old r=20 -> r=19); then add a "В ожидании" status to
# the second task row and bump the first task's progress percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was an empty (tall) spacer row; deleting it shifts the two data rows
# (old 3 -> 2, old 4 -> 3) up while preserving their styles/content/heights,
# and also shifts the trailing blank row (old 20 -> 19).
$ws.Rows.Item(2).Delete()

# The former row 20 (now row 19) was just a blank spacer row; remove it too.
$ws.Rows.Item(19).Delete()

# Bump the progress percentage on the first task row (formerly F3, now F2).
$ws.Range("F2").Value = 0.45

# Add a status note for the second task row (formerly row 4, now row 3).
$ws.Range("F3").Value = "В ожидании"

# Reflect the new active selection.
$ws.Range("C3").Select()
